$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select M1 (mirrors the authoring action of clicking into column M) then
# delete the whole column, which shifts column N left into M's place.
$ws.Range("M1").Select()
$ws.Columns.Item(13).Delete()
